$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1327.4286
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 1458.4
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1458.4
$ws.Range("M40").Value = -825
$ws.Range("N40").Value = -1808.4
$ws.Range("H43").Value = 3784
$ws.Range("I43").Value = 2980
$ws.Range("K43").Value = 2980
$ws.Range("M43").Value = -2911
$ws.Range("H53").Value = 532.35
$ws.Range("I53").Value = 213.45454
$ws.Range("J53").Value = 922.1111
$ws.Range("K53").Value = 213.45454
$ws.Range("L53").Value = 922.1111
$ws.Range("M53").Value = 423.54546
$ws.Range("N53").Value = -2196.1111
$ws.Range("H98").Value = 7085.45
$ws.Range("I98").Value = 4833.8887
$ws.Range("J98").Value = 8927.637000000001
$ws.Range("K98").Value = 4833.8887
$ws.Range("L98").Value = 8927.637000000001
$ws.Range("M98").Value = -3335.8887
$ws.Range("N98").Value = -11923.637
$ws.Range("H112").Value = 1275.3334
$ws.Range("J112").Value = 1275.3334
$ws.Range("L112").Value = 3826.0002
$ws.Range("N112").Value = -6042.0002
$ws.Range("H113").Value = 11228.3
$ws.Range("I113").Value = 6670.75
$ws.Range("J113").Value = 14266.667
$ws.Range("K113").Value = 6670.75
$ws.Range("L113").Value = 14266.667
$ws.Range("M113").Value = -3416.75
$ws.Range("N113").Value = -20774.667
$ws.Range("H122").Value = 7085.45
$ws.Range("I122").Value = 4833.8887
$ws.Range("J122").Value = 8927.637000000001
$ws.Range("K122").Value = 14501.6661
$ws.Range("L122").Value = 26782.911
$ws.Range("M122").Value = -12051.6661
$ws.Range("N122").Value = -31682.911
$ws.Range("H132").Value = 32710902
$ws.Range("I132").Value = 45633690
$ws.Range("K132").Value = 136901070
$ws.Range("M132").Value = -136898540
$ws.Range("H138").Value = 2602.7036
$ws.Range("I138").Value = 1706.909
$ws.Range("J138").Value = 3218.5625
$ws.Range("K138").Value = 5120.727000000001
$ws.Range("L138").Value = 9655.6875
$ws.Range("M138").Value = 19.27299999999923
$ws.Range("N138").Value = -19935.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 10755.25
$ws.Range("J43").Value = 10755.25
$ws.Range("L43").Value = 10755.25
$ws.Range("N43").Value = -11381.25
$ws.Range("H97").Value = 703.63635
$ws.Range("I97").Value = 637.7778
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 637.7778
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -141.7778
$ws.Range("N97").Value = -1992
$ws.Range("H133").Value = 37890
$ws.Range("J133").Value = 37890
$ws.Range("L133").Value = 37890
$ws.Range("N133").Value = -42950
$ws.Range("H135").Value = 95214.5
$ws.Range("J135").Value = 95214.5
$ws.Range("L135").Value = 95214.5
$ws.Range("N135").Value = -105354.5
$ws.Range("H137").Value = 48427.6
$ws.Range("J137").Value = 48427.6
$ws.Range("L137").Value = 48427.6
$ws.Range("N137").Value = -58627.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 8349778
$ws.Range("I7").Value = 17500000
$ws.Range("J7").Value = 5735429
$ws.Range("K7").Value = 17500000
$ws.Range("L7").Value = 5735429
$ws.Range("M7").Value = -17499887
$ws.Range("N7").Value = -5735655
$ws.Range("H99").Value = 7043.7896
$ws.Range("I99").Value = 1651.5
$ws.Range("J99").Value = 22142.2
$ws.Range("K99").Value = 1651.5
$ws.Range("L99").Value = 22142.2
$ws.Range("M99").Value = -153.5
$ws.Range("N99").Value = -25138.2
$ws.Range("H134").Value = 4894.706
$ws.Range("I134").Value = 2332
$ws.Range("J134").Value = 5962.5
$ws.Range("K134").Value = 6996
$ws.Range("L134").Value = 17887.5
$ws.Range("M134").Value = -4461
$ws.Range("N134").Value = -22957.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 57766.668
$ws.Range("J52").Value = 57766.668
$ws.Range("L52").Value = 57766.668
$ws.Range("N52").Value = -58354.668
$ws.Range("H134").Value = 2155.125
$ws.Range("I134").Value = 987.4545000000001
$ws.Range("J134").Value = 14999.5
$ws.Range("K134").Value = 2962.3635
$ws.Range("L134").Value = 44998.5
$ws.Range("M134").Value = -427.3635000000004
$ws.Range("N134").Value = -50068.5
$ws.Range("H137").Value = 43780
$ws.Range("J137").Value = 43780
$ws.Range("L137").Value = 43780
$ws.Range("N137").Value = -53980

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3014371.8
$ws.Range("I4").Value = 8036674.5
$ws.Range("J4").Value = 990
$ws.Range("K4").Value = 24110023.5
$ws.Range("L4").Value = 2970
$ws.Range("M4").Value = -24109911.5
$ws.Range("N4").Value = -3194
$ws.Range("H12").Value = 80.4375
$ws.Range("I12").Value = 10.666667
$ws.Range("J12").Value = 96.53846
$ws.Range("K12").Value = 32.000001
$ws.Range("L12").Value = 289.61538
$ws.Range("M12").Value = 140.999999
$ws.Range("N12").Value = -635.61538
$ws.Range("H22").Value = 579.8333
$ws.Range("I22").Value = 295.8
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 887.4000000000001
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = -718.4000000000001
$ws.Range("N22").Value = -6338
$ws.Range("H27").Value = 579.8333
$ws.Range("I27").Value = 295.8
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 887.4000000000001
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = -785.4000000000001
$ws.Range("N27").Value = -6204
$ws.Range("H33").Value = 140.11539
$ws.Range("I33").Value = 114.90909
$ws.Range("J33").Value = 158.6
$ws.Range("K33").Value = 689.4545400000001
$ws.Range("L33").Value = 951.5999999999999
$ws.Range("M33").Value = -406.4545400000001
$ws.Range("N33").Value = -1517.6
$ws.Range("H131").Value = 765.15
$ws.Range("J131").Value = 793.74725
$ws.Range("L131").Value = 2381.24175
$ws.Range("N131").Value = -12461.24175

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5923.9434
$ws.Range("I70").Value = 5576.8887
$ws.Range("J70").Value = 7876.125
$ws.Range("K70").Value = 5576.8887
$ws.Range("L70").Value = 7876.125
$ws.Range("M70").Value = -5306.8887
$ws.Range("N70").Value = -8416.125
$ws.Range("H73").Value = 5923.9434
$ws.Range("I73").Value = 5576.8887
$ws.Range("J73").Value = 7876.125
$ws.Range("K73").Value = 5576.8887
$ws.Range("L73").Value = 7876.125
$ws.Range("M73").Value = -4640.8887
$ws.Range("N73").Value = -9748.125
$ws.Range("H97").Value = 731.73334
$ws.Range("I97").Value = 633.36365
$ws.Range("J97").Value = 1002.25
$ws.Range("K97").Value = 633.36365
$ws.Range("L97").Value = 1002.25
$ws.Range("M97").Value = -137.36365
$ws.Range("N97").Value = -1994.25
$ws.Range("H126").Value = 3792.7144
$ws.Range("I126").Value = 2859.8
$ws.Range("J126").Value = 6125
$ws.Range("K126").Value = 8579.400000000001
$ws.Range("L126").Value = 18375
$ws.Range("M126").Value = -6109.400000000001
$ws.Range("N126").Value = -23315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5943.35
$ws.Range("I82").Value = 6414.8335
$ws.Range("J82").Value = 1700
$ws.Range("K82").Value = 6414.8335
$ws.Range("L82").Value = 1700
$ws.Range("M82").Value = -6053.8335
$ws.Range("N82").Value = -2422
$ws.Range("H85").Value = 5943.35
$ws.Range("I85").Value = 6414.8335
$ws.Range("J85").Value = 1700
$ws.Range("K85").Value = 6414.8335
$ws.Range("L85").Value = 1700
$ws.Range("M85").Value = -5166.8335
$ws.Range("N85").Value = -4196
$ws.Range("H100").Value = 5971.4287
$ws.Range("I100").Value = 2200
$ws.Range("J100").Value = 11000
$ws.Range("K100").Value = 2200
$ws.Range("L100").Value = 11000
$ws.Range("M100").Value = -1659
$ws.Range("N100").Value = -12082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11940
$ws.Range("J45").Value = 11625.5
$ws.Range("L45").Value = 11625.5
$ws.Range("N45").Value = -12607.5
$ws.Range("H100").Value = 1250.25
$ws.Range("I100").Value = 1000.3333
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 2000.6666
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -1459.6666
$ws.Range("N100").Value = -5082
